$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Model-year refresh: rows 2-9 and 53-54 move from 2020 -> 2021 and get
#    updated Base MSRP (column D) figures.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 2021
$ws.Range("D2").Value = 42120

$ws.Range("C3").Value = 2021
$ws.Range("D3").Value = 46590

$ws.Range("C4").Value = 2021
$ws.Range("D4").Value = 44810

$ws.Range("C5").Value = 2021
$ws.Range("D5").Value = 48765

$ws.Range("C6").Value = 2021
$ws.Range("D6").Value = 45050

$ws.Range("C7").Value = 2021
$ws.Range("D7").Value = 49520

$ws.Range("C8").Value = 2021
$ws.Range("D8").Value = 47215

$ws.Range("C9").Value = 2021
$ws.Range("D9").Value = 51130

$ws.Range("C53").Value = 2021
$ws.Range("D53").Value = 65875

$ws.Range("C54").Value = 2021
$ws.Range("D54").Value = 96675

# ---------------------------------------------------------------------------
# 2. Append four new "Black Line" trim rows (95-98) with their trim codes,
#    model years, MSRPs and DPHF fee - mirroring the layout/format of the
#    existing rows above them.
# ---------------------------------------------------------------------------
$currencyFormat = "$#,##0_);[Red]($#,##0)"
$dphfFormat = "$#,##0.00_);[Red]($#,##0.00)"

$newRows = @(
    @{ Row = 95; Code = "9203SE"; Trim = "RC 300 F SPORT Black Line";     Msrp = 48735 },
    @{ Row = 96; Code = "9207SE"; Trim = "RC 300 AWD F SPORT Black Line"; Msrp = 50910 },
    @{ Row = 97; Code = "9213SE"; Trim = "RC 350 F SPORT Black Line";     Msrp = 51665 },
    @{ Row = 98; Code = "9217SE"; Trim = "RC 350 AWD F SPORT Black Line"; Msrp = 53275 }
)

# New shared-string values are appended in the order they are first written,
# so populate every trim-code cell (column A) before any trim-name cell
# (column B) to reproduce the original author's string table ordering.
foreach ($r in $newRows) {
    $ws.Range("A$($r.Row)").Value = $r.Code
}
foreach ($r in $newRows) {
    $ws.Range("B$($r.Row)").Value = $r.Trim
}
foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("C$row").Value = 2021
    $ws.Range("D$row").Value = $r.Msrp
    $ws.Range("D$row").NumberFormat = $currencyFormat
    $ws.Range("E$row").Value = 1025
    $ws.Range("E$row").NumberFormat = $dphfFormat
}

# ---------------------------------------------------------------------------
# 3. Update the sheet view state (scroll position / active selection) to
#    reflect where the author was working while testing the new rows.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C55").Select()
